# 补做 PPO-6 与纠正 DPPO-K5-early stop 的实验
# Adds two new rows (12 and 13) of experiment-log data to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 12 - corrected DPPO K5 early-stop run
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "result/DPPO/DPPO_warm_1000_run_300_simulate_1024_plan_5_stop_l_10_rate_0.9_correct"
# Copy the "green" result-path font/style used elsewhere in column A (e.g. A2)
# onto A12 without disturbing the value we just wrote.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

$ws.Range("B12").Value = "DPPO"
$ws.Range("C12").Value = 300
$ws.Range("D12").Value = 1000
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 0.5
$ws.Range("G12").Value = 3000
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = "是"
$ws.Range("J12").Value = "early stop 的正确实验。但是都是在31轮停止的world model"

# ---------------------------------------------------------------------
# Row 13 - PPO x6 data-volume comparison run
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "result/PPO/PPO_warm_1000_run_500_simulate_6144"
# Start from the existing "blue-grey" result-path font (style used by A9/A10),
# then recolour it to a grey tone - this yields a brand-new font/style entry,
# matching the newly introduced font in the workbook.
$ws.Range("A9").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Font.Color = 8421504

$ws.Range("B13").Value = "PPO"
$ws.Range("C13").Value = 500
$ws.Range("D13").Value = 1000
$ws.Range("J13").Value = "每个eisode使用的数据量是6倍（6144/1024=6）"

# Final selection lands on J13, matching where the author finished editing.
$ws.Range("J13").Select() | Out-Null
